$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 & 5: swap company names ---
$ws.Range("B4").Value = "C Security Systems AB (publ) (NGM:CSEC)"
$ws.Range("B5").Value = "MultiDocker Cargo Handling AB (publ) (NGM:MULT)"

# --- Row 5: remove D5 (historical_growth_revenue_last_5_years) ---
$ws.Range("D5").ClearContents()

# --- Row 2 & 3: remove T (buybacks_cash_returned) data cells ---
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()

# --- Row 2 updated values ---
$ws.Range("D2").Value = -0.02650000000000001
$ws.Range("G2").Value = 0.1311455108359133
$ws.Range("H2").Value = 0.1191331269349845
$ws.Range("I2").Value = -0.3497213622291021
$ws.Range("J2").Value = -0.3497213622291021
$ws.Range("K2").Value = -13.823
$ws.Range("L2").Value = -0.2853044375644994
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 17.192
$ws.Range("V2").Value = 0.252515312192471
$ws.Range("W2").Value = -1.076712328767123
$ws.Range("X2").Value = 0.05053975211698089
$ws.Range("Y2").Value = -1.127252080884104
$ws.Range("Z2").Value = 0.2493720662109858
$ws.Range("AA2").Value = -1.145945945945946
$ws.Range("AB2").Value = 0.04634885448143908
$ws.Range("AC2").Value = -1.192294800427385
$ws.Range("AD2").Value = 3.076
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3.076
$ws.Range("AG2").Value = -14.116
$ws.Range("AH2").Value = 0.04322713922342923
$ws.Range("AI2").Value = 0.01481174731428737
$ws.Range("AJ2").Value = -0.261567254062668
$ws.Range("AK2").Value = -0.07410712879499792
$ws.Range("AL2").Value = 0.788
$ws.Range("AM2").Value = 0.4710000000000001
$ws.Range("AN2").Value = -0.3910998092816275
$ws.Range("AO2").Value = -21.50253807106599
$ws.Range("AP2").Value = 1.794787031150668
$ws.Range("AQ2").Value = -35.97452229299362

# --- Row 3 updated values ---
$ws.Range("D3").Value = -0.25
$ws.Range("G3").Value = 0.2714617169373549
$ws.Range("H3").Value = 0.2714617169373549
$ws.Range("I3").Value = -0.234338747099768
$ws.Range("J3").Value = -0.234338747099768
$ws.Range("K3").Value = -10.2
$ws.Range("L3").Value = -0.2366589327146171
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 17
$ws.Range("V3").Value = 0.2627511591962906
$ws.Range("W3").Value = -0.04843304843304843
$ws.Range("X3").Value = 0.04557491649744751
$ws.Range("Y3").Value = -0.09400796493049593
$ws.Range("Z3").Value = 0.2264797376828653
$ws.Range("AA3").Value = -0.05307297797208677
$ws.Range("AB3").Value = 0.04535861760680459
$ws.Range("AC3").Value = -0.09843159557889136
$ws.Range("AD3").Value = 0.5580000000000001
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.5580000000000001
$ws.Range("AG3").Value = -16.442
$ws.Range("AH3").Value = 0.008550675779214808
$ws.Range("AI3").Value = 0.002749337301313573
$ws.Range("AJ3").Value = -0.3407103485432467
$ws.Range("AK3").Value = -0.08841781477537938
$ws.Range("AL3").Value = 0.34
$ws.Range("AM3").Value = 0.05600000000000005
$ws.Range("AN3").Value = -0.4810344827586208
$ws.Range("AO3").Value = -29.70588235294117
$ws.Range("AP3").Value = 14.17413793103448
$ws.Range("AQ3").Value = -180.3571428571427

# --- Row 4 updated values ---
$ws.Range("D4").Value = 0.197
$ws.Range("G4").Value = -1.911111111111111
$ws.Range("H4").Value = -2.044444444444444
$ws.Range("I4").Value = -2.355555555555556
$ws.Range("J4").Value = -2.355555555555556
$ws.Range("K4").Value = -0.393
$ws.Range("L4").Value = -2.183333333333334
$ws.Range("U4").Value = 0.022
$ws.Range("V4").Value = 0.008835341365461847
$ws.Range("W4").Value = -1.076712328767123
$ws.Range("X4").Value = 0.05053975211698089
$ws.Range("Y4").Value = -1.127252080884104
$ws.Range("Z4").Value = 0.4864864864864865
$ws.Range("AA4").Value = -1.145945945945946
$ws.Range("AB4").Value = 0.04634885448143908
$ws.Range("AC4").Value = -1.192294800427385
$ws.Range("AD4").Value = 0.458
$ws.Range("AF4").Value = 0.458
$ws.Range("AG4").Value = 0.436
$ws.Range("AH4").Value = 0.155359565807327
$ws.Range("AI4").Value = 0.1925988225399495
$ws.Range("AJ4").Value = 0.1490088858509911
$ws.Range("AK4").Value = 0.1850594227504244
$ws.Range("AL4").Value = 0.013
$ws.Range("AM4").Value = -0.02
$ws.Range("AN4").Value = -1.159493670886076
$ws.Range("AO4").Value = -32.61538461538461
$ws.Range("AP4").Value = -1.10379746835443
$ws.Range("AQ4").Value = 21.2

# --- Row 5 updated values ---
$ws.Range("G5").Value = -0.9675048355899419
$ws.Range("H5").Value = -1.075435203094778
$ws.Range("I5").Value = -1.241779497098646
$ws.Range("J5").Value = -1.241779497098646
$ws.Range("K5").Value = -3.23
$ws.Range("L5").Value = -0.6247582205029013
$ws.Range("U5").Value = 0.17
$ws.Range("V5").Value = 0.1903695408734603
$ws.Range("W5").Value = -2.543307086614173
$ws.Range("X5").Value = 0.1106603728085735
$ws.Range("Y5").Value = -2.653967459422747
$ws.Range("Z5").Value = 1.43054786939679
$ws.Range("AA5").Value = -1.776425013835086
$ws.Range("AB5").Value = 0.04761057938301933
$ws.Range("AC5").Value = -1.824035593218105
$ws.Range("AD5").Value = 2.06
$ws.Range("AF5").Value = 2.06
$ws.Range("AG5").Value = 1.89
$ws.Range("AH5").Value = 0.6975956654249915
$ws.Range("AI5").Value = 0.881471972614463
$ws.Range("AJ5").Value = 0.6791232482932087
$ws.Range("AK5").Value = 0.8721735117674203
$ws.Range("AL5").Value = 0.435
$ws.Range("AM5").Value = 0.435
$ws.Range("AN5").Value = -0.3264659270998416
$ws.Range("AO5").Value = -14.75862068965517
$ws.Range("AP5").Value = -0.2995245641838352
$ws.Range("AQ5").Value = -14.75862068965517

